$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date labels for the 6 newly appended rows (214-219)
$dates = @("26-10-2021", "27-10-2021", "28-10-2021", "29-10-2021", "01-11-2021", "02-11-2021")
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item(214 + $i, 1).Value = $dates[$i]
}

# Data rows, keyed by row number -> column letter -> value.
# Row 213 gets additional (previously-blank) values filled in; rows
# 214-219 are brand-new rows continuing the daily series.
$rowData = @{
    213 = @{ B=0.25; C=0.1; D=-0.1; E=-0.5; F=0.75; G=4.35; H=1.75; I=0.5; J=1.5; K=7.5; L=0.5; M=1.13; N=16; O=38; P=6.25; Q=2;   R=4.75; S=1.5 }
    214 = @{ B=0.25; C=0.1; D=-0.1; E=-0.5; F=0.75; G=4.35; H=1.75; I=0.5; J=1.5; K=7.5; L=0.5; M=1.13; N=16; O=38; P=6.25; Q=2;   R=4.75; S=1.5 }
    215 = @{ B=0.25; C=0.1; D=-0.1; E=-0.5; F=0.75; G=4.35; H=1.75; I=0.5; J=1.5; K=7.5; L=0.5; M=1.13; N=16; O=38; P=7.75; Q=2;   R=4.75; S=1.5 }
    216 = @{ B=0.25; C=0.1; D=-0.1; E=-0.5; F=0.75; G=4.35; H=1.75;               K=7.5; L=0.5; M=1.13; N=16; O=38; P=7.75; Q=2;   R=4.75; S=1.5 }
    217 = @{ B=0.25; C=0.1; D=-0.1; E=-0.5; F=0.75; G=4.35; H=1.75; I=0.5; J=1.5; K=7.5; L=0.5; M=1.13; N=16; O=38; P=7.75; Q=2.5; R=4.75; S=1.5 }
    218 = @{ B=0.25; C=0.1; D=-0.1; E=-0.5; F=0.75; G=4.35; H=1.75; I=0.5; J=1.5; K=7.5; L=0.5; M=1.13; N=16; O=38; P=7.75; Q=2.5; R=4.75 }
    219 = @{ B=0.25; C=0.1;                         F=0.75; G=4.35; H=1.75; I=0.5; J=1.5; K=7.5; L=0.5;               N=16;                   R=4.75; S=1.5 }
}

$colIndex = @{ A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13; N=14; O=15; P=16; Q=17; R=18; S=19 }

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($colLetter in $cols.Keys) {
        $c = $colIndex[$colLetter]
        $ws.Cells.Item($r, $c).Value = $cols[$colLetter]
    }
}
